# A new weekly price record was inserted at the top of the data table
# (row 13, right after the header block). Every existing record from the
# old row 13 onward shifts down by one row, and the brand-new record is
# written into the vacated row 13.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push row 13 and everything below it down by one row, opening up a
# blank row 13 for the new record.
$ws.Rows.Item(13).Insert()

# The row that used to be 13 is now row 14; duplicate it into the new
# row 13 so all the "carried over" columns (Mercado, Región, Codreg,
# Categoría, Variedad, Calidad, Precio mínimo/máximo, Unidad, Origen,
# Kg o Unidades, Clasificación) start out correct.
$ws.Range("A14:R14").Copy()
$ws.Range("A13:R13").PasteSpecial()

# Now overwrite the fields that actually differ for the new record:
# date, volume, weighted average price and price per kg.
$ws.Range("D13").Value = 45111
$ws.Range("J13").Value = 350
$ws.Range("M13").Value = 1471
$ws.Range("P13").Value = 490
